# Apply edits described in commit: "Adjustment to prevent additional coal
# CCS retrofits beyond EPA projections"
#
# 1. On the "About" sheet, append three new note rows (14-16) explaining
#    why no additional coal CCS retrofits are allowed.
# 2. On the "BBPPRTY" sheet, the "hard coal w CCS" row (row 19) is changed
#    so every year (columns B:AE) is set to 1 (previously it switched to 0
#    starting in 2028/column I), and the cells' existing number format is
#    cleared so they revert to the workbook's default (unformatted) style.

$wb = $excel.ActiveWorkbook

# --- "About" sheet: add explanatory notes -----------------------------
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A14").Value = "We already include all coal CCS retrofits projected by EPA as part of its 111 rules."
$wsAbout.Range("A15").Value = "Projected coal retriements and retrofits cover all of the exisiting coal fleet between"
$wsAbout.Range("A16").Value = "2028 and 2035, so we do not allow any additional retrofits."

# --- "BBPPRTY" sheet: update "hard coal w CCS" retrofit-ban row -------
$wsData = $wb.Worksheets.Item("BBPPRTY")

$rowRange = $wsData.Range("B19:AE19")
$rowRange.ClearFormats()
$rowRange.Value = 1

# Update the cell selections on each sheet to mirror the authored edit,
# then leave "About" as the active sheet/tab (as in the source workbook).
$wsData.Range("B29").Select()

$wsAbout.Range("A16").Select()
$wsAbout.Activate()
